$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.991.79"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "3.560.57"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'190.95"
$ws.Range("E5").Value = "  -0.97%  "

$ws.Range("D6").Value = "'568.44"
$ws.Range("E6").Value = "  -5.73%  "

$ws.Range("D7").Value = "3.556.23"
$ws.Range("E7").Value = "  -2.99%  "

$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "'0.674"

$ws.Range("D11").Value = "'55.66"
$ws.Range("E11").Value = "  -3.75%  "

$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").Value = "'9.82"
$ws.Range("E14").Value = "  -3.77%  "

$ws.Range("D15").Value = "4.128.09"
$ws.Range("E15").Value = "  -3.19%  "

$ws.Range("D16").Value = "3.566.17"
$ws.Range("E16").Value = "  -2.94%  "

$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").Value = "66.902.32"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("D19").Value = "'12.13"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("D20").Value = "'18.13"
$ws.Range("E20").Value = "  -4.33%  "

$ws.Range("E21").Value = "  -5.43%  "

$ws.Range("D22").Value = "'399.56"
$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("E23").Value = "  -7.21%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "'11.82"
$ws.Range("E24").Value = "  +5.23%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'85.52"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("D26").Value = "'2.90"
$ws.Range("E26").Value = "  -2.05%  "

$ws.Range("D27").Value = "'12.42"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").Value = "'3.68"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "'7.76"
$ws.Range("E30").Value = "  +4.26%  "

$ws.Range("D31").Value = "'8.95"
$ws.Range("E31").Value = "  -3.72%  "

$ws.Range("D32").Value = "'31.08"
$ws.Range("E32").Value = "  -2.62%  "

$ws.Range("D33").Value = "'641.26"
$ws.Range("E33").Value = "  +5.24%  "

$ws.Range("D34").Value = "'12.09"
$ws.Range("E34").Value = "  -2.15%  "

$ws.Range("E35").Value = "  -2.96%  "

$ws.Range("D36").Value = "'63.67"
$ws.Range("E36").Value = "  -6.77%  "

$ws.Range("D37").Value = "'42.11"
$ws.Range("E37").Value = "  -6.83%  "

$ws.Range("D38").Value = "'0.403"
$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "0.0₃0759"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").Value = "3.171.46"
$ws.Range("E41").Value = "  +12.45%  "

$ws.Range("E42").Value = "  -1.37%  "

$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "'2.68"
$ws.Range("E45").Value = "  +5.01%  "

$ws.Range("D46").Value = "'0.0412"
$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.11"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.130"
$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("D49").Value = "'141.82"
$ws.Range("E49").Value = "  -1.86%  "

$ws.Range("D50").Value = "'8.49"
$ws.Range("E50").Value = "  -5.52%  "

$ws.Range("D51").Value = "'2.52"
$ws.Range("E51").Value = "  -4.25%  "
